$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '38.760.03'
$ws.Range("E2").Value = '  +0.36%  '

$ws.Range("D3").Value = '2.102.21'
$ws.Range("E3").Value = '  +0.45%  '

$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '227.70'
$ws.Range("E5").Value = '  -0.52%  '

$ws.Range("E6").Value = '  +0.41%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '62.40'
$ws.Range("E7").Value = '  +1.69%  '

$ws.Range("E8").Value = '  -0.07%  '

$ws.Range("E9").Value = '  +2.25%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0841'
$ws.Range("E10").Value = '  -0.20%  '

$ws.Range("E11").Value = '  -1.04%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.74'
$ws.Range("E12").Value = '  +6.42%  '

$ws.Range("D13").Value = '2.413.55'
$ws.Range("E13").Value = '  +0.46%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.02'
$ws.Range("E14").Value = '  -1.33%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.808'
$ws.Range("E15").Value = '  +3.22%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.53'
$ws.Range("E16").Value = '  +1.35%  '

$ws.Range("D17").Value = '2.103.70'
$ws.Range("E17").Value = '  -0.29%  '

$ws.Range("D18").Value = '38.746.77'
$ws.Range("E18").Value = '  +0.49%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '71.83'
$ws.Range("E19").Value = '  +1.13%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.12'
$ws.Range("E20").Value = '  +0.99%  '

$ws.Range("E21").Value = '  +0.51%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '227.66'
$ws.Range("E22").Value = '  +0.67%  '

$ws.Range("E24").Value = '  -4.06%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.32'
$ws.Range("E25").Value = '  -0.33%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.63'
$ws.Range("E26").Value = '  +2.03%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '172.18'
$ws.Range("E27").Value = '  +1.01%  '

$ws.Range("E28").Value = '  +3.95%  '

$ws.Range("E29").Value = '  +4.29%  '

$ws.Range("E30").Value = '  +1.19%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.51'
$ws.Range("E31").Value = '  +10.17%  '

$ws.Range("E32").Value = '  +0.52%  '

$ws.Range("E33").Value = '  +1.20%  '

$ws.Range("E34").Value = '  -0.92%  '

$ws.Range("E35").Value = '  +7.02%  '

$ws.Range("E36").Value = '  +1.80%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.40'
$ws.Range("E37").Value = '  +0.50%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.58'
$ws.Range("E38").Value = '  -0.23%  '

$ws.Range("E39").Value = '  -0.07%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.06'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '102.86'
$ws.Range("E41").Value = '  +3.05%  '

$ws.Range("E42").Value = '  +4.10%  '

$ws.Range("D43").Value = '1.528.43'
$ws.Range("E43").Value = '  -1.11%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.21'
$ws.Range("E44").Value = '  +7.77%  '

$ws.Range("E45").Value = '  -1.04%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.81'
$ws.Range("E46").Value = '  +1.65%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0913'
$ws.Range("E47").Value = '  +0.15%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.17'
$ws.Range("E48").Value = '  +0.64%  '

$ws.Range("E49").Value = '  +1.90%  '

$ws.Range("E50").Value = '  -0.85%  '

$ws.Range("D51").Value = '2.300.07'
$ws.Range("E51").Value = '  +0.37%  '
